$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: Locality corrected from "St. Kilda" to "Foula" ---
$ws.Range("C17").Value = "Foula"

# --- Row 20: year range corrected from "2005-2006" to "2001-2002" ---
$ws.Range("B20").Value = "2001-2002"

# --- Row 21 (new): Hyperoplus lanceolatus, North Sea record ---
$ws.Range("A21").Value = "Hyperoplus lanceolatus"
$ws.Range("B21").Value = "2001-2001"
$ws.Range("C21").Value = "North Sea"
$ws.Range("D21").Value = 7
$ws.Range("E21").Value = 16.1
$ws.Range("F21").Value = 1.3
$ws.Range("G21").Value = -16.4
$ws.Range("H21").Value = 0.4
$ws.Range("I21").Value = "Das et al., 2003 https://www.int-res.com/articles/meps2003/263/m263p287.pdf"

# --- Row 22 (new): Ammodytes marinus, Foula, Bearhop et al., 1999 ---
$ws.Range("A22").Value = "Ammodytes marinus"
$ws.Range("B22").Value = 1996
$ws.Range("C22").Value = "Foula"
$ws.Range("D22").Value = 12
$ws.Range("E22").Value = 7.9
$ws.Range("F22").Value = 0.95
$ws.Range("G22").Value = -17.5
$ws.Range("H22").Value = 0.75
$ws.Range("I22").Value = "Bearhop et al., 1999 https://www.jstor.org/stable/2655696?casa_token=LG2WgTk6tNcAAAAA%3AiksV4fDbxPgsqLOjlGVdUJJ0P7bQT2GQLaNpNMnLGNnCJwTF3AXbc5Txl_jvD7BSSC0uNJN28d9C6WHIQkFUHMP11eYahJprlNIsXvtzwes6bS61Ig&seq=4"

# Move the active selection to reflect the end of the edit session
$ws.Range("F24").Select()
